$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 58: add values in columns C, D, E (E58 first so new shared string ordering matches)
$ws.Range("E58").Value = "what is this? never shows a value in any of my inputs"
$ws.Range("C58").Value = "?"
$ws.Range("D58").Value = "?"

# Row 57: add comment in column E
$ws.Range("E57").Value = "needed? very simple test"

# Update sheet view: top left cell and selection
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E58").Select()
